$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-price record was added for this market/product. It belongs
# chronologically before the existing row 339, so insert a fresh row there
# and let Excel shift rows 339-409 down to 340-410 (preserving formats,
# including the date NumberFormat on column D).
$ws.Rows(339).Insert()

$ws.Range("A339").Value = 4
$ws.Range("B339").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C339").Value = "Los Lagos"
$ws.Range("D339").Value = 45204
$ws.Range("E339").Value = 10
$ws.Range("F339").Value = 100112039
$ws.Range("G339").Value = "Ciboulette"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 80
$ws.Range("K339").Value = 3500
$ws.Range("L339").Value = 3500
$ws.Range("M339").Value = 3500
$ws.Range("N339").Value = "`$/docena de atados"
$ws.Range("O339").Value = "Región Metropolitana"
$ws.Range("P339").Value = 1167
$ws.Range("Q339").Value = 3
$ws.Range("R339").Value = "Hortaliza"
